$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2391.1667
$ws.Range("I2").Value = 2469.6
$ws.Range("J2").Value = 1999
$ws.Range("K2").Value = 2469.6
$ws.Range("L2").Value = 1999
$ws.Range("M2").Value = -2356.6
$ws.Range("N2").Value = -2225
$ws.Range("H15").Value = 689.3273
$ws.Range("I15").Value = 689.3273
$ws.Range("K15").Value = 2067.9819
$ws.Range("M15").Value = -1898.9819
$ws.Range("H29").Value = 3630.6
$ws.Range("I29").Value = 76.5
$ws.Range("J29").Value = 6000
$ws.Range("K29").Value = 229.5
$ws.Range("L29").Value = 18000
$ws.Range("M29").Value = 51.5
$ws.Range("N29").Value = -18562
$ws.Range("H38").Value = 61.25
$ws.Range("I38").Value = 61.25
$ws.Range("K38").Value = 183.75
$ws.Range("M38").Value = 188.25
$ws.Range("H43").Value = 2050.8333
$ws.Range("I43").Value = 1301
$ws.Range("K43").Value = 1301
$ws.Range("M43").Value = -1232
$ws.Range("H58").Value = 3167
$ws.Range("I58").Value = 2483.8
$ws.Range("K58").Value = 7451.400000000001
$ws.Range("M58").Value = -7301.400000000001
$ws.Range("H86").Value = 2199.6667
$ws.Range("I86").Value = 675
$ws.Range("J86").Value = 5249
$ws.Range("K86").Value = 675
$ws.Range("L86").Value = 5249
$ws.Range("M86").Value = 448
$ws.Range("N86").Value = -7495
$ws.Range("H89").Value = 2199.6667
$ws.Range("I89").Value = 675
$ws.Range("J89").Value = 5249
$ws.Range("K89").Value = 3375
$ws.Range("L89").Value = 26245
$ws.Range("M89").Value = 2241
$ws.Range("N89").Value = -37477
$ws.Range("H98").Value = 1184.5405
$ws.Range("I98").Value = 1184.5405
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1184.5405
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 313.4594999999999
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 1184.5405
$ws.Range("I122").Value = 1184.5405
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3553.6215
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1103.6215
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 37040076
$ws.Range("I132").Value = 45457140
$ws.Range("J132").Value = 4991.2
$ws.Range("K132").Value = 136371420
$ws.Range("L132").Value = 14973.6
$ws.Range("M132").Value = -136368890
$ws.Range("N132").Value = -20033.6
$ws.Range("H137").Value = 138641.69
$ws.Range("I137").Value = 223620
$ws.Range("K137").Value = 670860
$ws.Range("M137").Value = -668310
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3477.9167
$ws.Range("J32").Value = 5570
$ws.Range("L32").Value = 5570
$ws.Range("N32").Value = -6144
$ws.Range("H61").Value = 3608.303
$ws.Range("I61").Value = 3564.9688
$ws.Range("K61").Value = 3564.9688
$ws.Range("M61").Value = -3352.9688
$ws.Range("H122").Value = 3475682.2
$ws.Range("I122").Value = 3947.5
$ws.Range("K122").Value = 11842.5
$ws.Range("M122").Value = -9392.5
$ws.Range("H132").Value = 3302.4167
$ws.Range("I132").Value = 2375.5908
$ws.Range("K132").Value = 7126.7724
$ws.Range("M132").Value = -4596.7724
$ws.Range("H136").Value = 3608.303
$ws.Range("I136").Value = 3564.9688
$ws.Range("K136").Value = 10694.9064
$ws.Range("M136").Value = -8144.9064
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13428.161
$ws.Range("I31").Value = 1291.25
$ws.Range("J31").Value = 16500.797
$ws.Range("K31").Value = 1291.25
$ws.Range("L31").Value = 16500.797
$ws.Range("M31").Value = -996.25
$ws.Range("N31").Value = -17090.797
$ws.Range("H34").Value = 13428.161
$ws.Range("I34").Value = 1291.25
$ws.Range("J34").Value = 16500.797
$ws.Range("K34").Value = 1291.25
$ws.Range("L34").Value = 16500.797
$ws.Range("M34").Value = -1089.25
$ws.Range("N34").Value = -16904.797
$ws.Range("H99").Value = 5954.091
$ws.Range("I99").Value = 5857.143
$ws.Range("J99").Value = 6123.75
$ws.Range("K99").Value = 5857.143
$ws.Range("L99").Value = 6123.75
$ws.Range("M99").Value = -4359.143
$ws.Range("N99").Value = -9119.75
$ws.Range("H103").Value = 42988.223
$ws.Range("I103").Value = 42988.223
$ws.Range("K103").Value = 42988.223
$ws.Range("M103").Value = -41816.223
$ws.Range("H122").Value = 3807.0625
$ws.Range("J122").Value = 4619.4443
$ws.Range("L122").Value = 13858.3329
$ws.Range("N122").Value = -18758.3329
$ws.Range("H126").Value = 5954.091
$ws.Range("I126").Value = 5857.143
$ws.Range("J126").Value = 6123.75
$ws.Range("K126").Value = 17571.429
$ws.Range("L126").Value = 18371.25
$ws.Range("M126").Value = -15101.429
$ws.Range("N126").Value = -23311.25
$ws.Range("H132").Value = 74383.71000000001
$ws.Range("I132").Value = 79413.234
$ws.Range("K132").Value = 238239.702
$ws.Range("M132").Value = -235709.702
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 782.25
$ws.Range("I8").Value = 782.25
$ws.Range("K8").Value = 2346.75
$ws.Range("M8").Value = -2207.75
$ws.Range("H12").Value = 129417.43
$ws.Range("I12").Value = 296301
$ws.Range("K12").Value = 888903
$ws.Range("M12").Value = -888730
$ws.Range("H68").Value = 2192.8572
$ws.Range("J68").Value = 2800
$ws.Range("L68").Value = 8400
$ws.Range("N68").Value = -10022
$ws.Range("H71").Value = 2192.8572
$ws.Range("J71").Value = 2800
$ws.Range("L71").Value = 25200
$ws.Range("N71").Value = -33312
$ws.Range("H94").Value = 7378.143
$ws.Range("H107").Value = 1228.4
$ws.Range("J107").Value = 1435
$ws.Range("L107").Value = 4305
$ws.Range("N107").Value = -8145
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 91465910
$ws.Range("I80").Value = 109758600
$ws.Range("J80").Value = 2506
$ws.Range("K80").Value = 109758600
$ws.Range("L80").Value = 2506
$ws.Range("M80").Value = -109757602
$ws.Range("N80").Value = -4502
$ws.Range("H83").Value = 91465910
$ws.Range("I83").Value = 109758600
$ws.Range("J83").Value = 2506
$ws.Range("K83").Value = 548793000
$ws.Range("L83").Value = 12530
$ws.Range("M83").Value = -548788008
$ws.Range("N83").Value = -22514
$ws.Range("H102").Value = 4857050
$ws.Range("I102").Value = 8549154
$ws.Range("J102").Value = 1657225.8
$ws.Range("K102").Value = 8549154
$ws.Range("L102").Value = 1657225.8
$ws.Range("M102").Value = -8547532
$ws.Range("N102").Value = -1660469.8
$ws.Range("H122").Value = 2964648.2
$ws.Range("I122").Value = 2964648.2
$ws.Range("K122").Value = 8893944.600000001
$ws.Range("M122").Value = -8891494.600000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1606.5714
$ws.Range("I16").Value = 1606.5714
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1606.5714
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1436.5714
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value = 11499.929
$ws.Range("I40").Value = 9589
$ws.Range("K40").Value = 9589
$ws.Range("M40").Value = -9453
$ws.Range("H93").Value = 17544880
$ws.Range("I93").Value = 25642012
$ws.Range("J93").Value = 1098.3334
$ws.Range("K93").Value = 25642012
$ws.Range("L93").Value = 1098.3334
$ws.Range("M93").Value = -25640764
$ws.Range("N93").Value = -3594.3334
$ws.Range("H100").Value = 5529.2354
$ws.Range("I100").Value = 4999.769
$ws.Range("K100").Value = 4999.769
$ws.Range("M100").Value = -4458.769
$ws.Range("H122").Value = 7975.3335
$ws.Range("I122").Value = 4998.3335
$ws.Range("J122").Value = 9463.833000000001
$ws.Range("K122").Value = 14995.0005
$ws.Range("L122").Value = 28391.499
$ws.Range("M122").Value = -12545.0005
$ws.Range("N122").Value = -33291.499
$ws.Range("H136").Value = 27011.477
$ws.Range("I136").Value = 31187.486
$ws.Range("K136").Value = 93562.458
$ws.Range("M136").Value = -91012.458
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H96").Value = 2000
$ws.Range("I96").Value = 2000
$ws.Range("K96").Value = 2000
$ws.Range("M96").Value = -627
$ws.Range("H122").Value = 2777.8572
$ws.Range("I122").Value = 2625.9678
$ws.Range("K122").Value = 7877.903399999999
$ws.Range("M122").Value = -5427.903399999999
$ws.Range("H126").Value = 3419.2856
$ws.Range("I126").Value = 2666
$ws.Range("J126").Value = 5302.5
$ws.Range("K126").Value = 7998
$ws.Range("L126").Value = 15907.5
$ws.Range("M126").Value = -5528
$ws.Range("N126").Value = -20847.5
$ws.Range("H136").Value = 3693.1667
$ws.Range("I136").Value = 3683.535
$ws.Range("J136").Value = 3730.818
$ws.Range("K136").Value = 11050.605
$ws.Range("L136").Value = 11192.454
$ws.Range("M136").Value = -8500.605
$ws.Range("N136").Value = -16292.454
